$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A 2980-2026
$ws.Range("A2").Value = 'A 2980-2026'
$ws.Range("B2").Value = 46038.62965277778
$ws.Range("C2").Value = 46074
$ws.Range("D2").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E2").Value = 'SOTENÄS'
$ws.Range("G2").Value = 1.1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 'Blåsippa'
$ws.Range("S2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/artfynd/A 2980-2026 artfynd.xlsx`", `"A 2980-2026`")"
$ws.Range("T2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/kartor/A 2980-2026 karta.png`", `"A 2980-2026`")"
$ws.Range("V2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomål/A 2980-2026 FSC-klagomål.docx`", `"A 2980-2026`")"
$ws.Range("W2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomålsmail/A 2980-2026 FSC-klagomål mail.docx`", `"A 2980-2026`")"
$ws.Range("X2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsyn/A 2980-2026 tillsynsbegäran.docx`", `"A 2980-2026`")"
$ws.Range("Y2").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsynsmail/A 2980-2026 tillsynsbegäran mail.docx`", `"A 2980-2026`")"

# Row 3: A 24087-2025
$ws.Range("A3").Value = 'A 24087-2025'
$ws.Range("B3").Value = 45795
$ws.Range("C3").Value = 46074
$ws.Range("D3").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E3").Value = 'SOTENÄS'
$ws.Range("G3").Value = 30.3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 'Rödlånke'
$ws.Range("S3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/artfynd/A 24087-2025 artfynd.xlsx`", `"A 24087-2025`")"
$ws.Range("T3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/kartor/A 24087-2025 karta.png`", `"A 24087-2025`")"
$ws.Range("V3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomål/A 24087-2025 FSC-klagomål.docx`", `"A 24087-2025`")"
$ws.Range("W3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomålsmail/A 24087-2025 FSC-klagomål mail.docx`", `"A 24087-2025`")"
$ws.Range("X3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsyn/A 24087-2025 tillsynsbegäran.docx`", `"A 24087-2025`")"
$ws.Range("Y3").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsynsmail/A 24087-2025 tillsynsbegäran mail.docx`", `"A 24087-2025`")"

# Row 4: A 24088-2025
$ws.Range("A4").Value = 'A 24088-2025'
$ws.Range("B4").Value = 45795
$ws.Range("C4").Value = 46074
$ws.Range("D4").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E4").Value = 'SOTENÄS'
$ws.Range("G4").Value = 2.8
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 'Blåsippa'
$ws.Range("S4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/artfynd/A 24088-2025 artfynd.xlsx`", `"A 24088-2025`")"
$ws.Range("T4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/kartor/A 24088-2025 karta.png`", `"A 24088-2025`")"
$ws.Range("V4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomål/A 24088-2025 FSC-klagomål.docx`", `"A 24088-2025`")"
$ws.Range("W4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomålsmail/A 24088-2025 FSC-klagomål mail.docx`", `"A 24088-2025`")"
$ws.Range("X4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsyn/A 24088-2025 tillsynsbegäran.docx`", `"A 24088-2025`")"
$ws.Range("Y4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsynsmail/A 24088-2025 tillsynsbegäran mail.docx`", `"A 24088-2025`")"

# Row 5: A 24036-2025
$ws.Range("A5").Value = 'A 24036-2025'
$ws.Range("B5").Value = 45795
$ws.Range("C5").Value = 46074
$ws.Range("D5").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E5").Value = 'SOTENÄS'
$ws.Range("G5").Value = 37.1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 'Kopparödla'
$ws.Range("S5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/artfynd/A 24036-2025 artfynd.xlsx`", `"A 24036-2025`")"
$ws.Range("T5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/kartor/A 24036-2025 karta.png`", `"A 24036-2025`")"
$ws.Range("V5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomål/A 24036-2025 FSC-klagomål.docx`", `"A 24036-2025`")"
$ws.Range("W5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/klagomålsmail/A 24036-2025 FSC-klagomål mail.docx`", `"A 24036-2025`")"
$ws.Range("X5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsyn/A 24036-2025 tillsynsbegäran.docx`", `"A 24036-2025`")"
$ws.Range("Y5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1427/tillsynsmail/A 24036-2025 tillsynsbegäran mail.docx`", `"A 24036-2025`")"

# Row 6: A 18791-2021
$ws.Range("A6").Value = 'A 18791-2021'
$ws.Range("B6").Value = 44307.58013888889
$ws.Range("C6").Value = 46074
$ws.Range("D6").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E6").Value = 'SOTENÄS'
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0

# Row 7: A 54967-2022
$ws.Range("A7").Value = 'A 54967-2022'
$ws.Range("B7").Value = 44886.4831712963
$ws.Range("C7").Value = 46074
$ws.Range("D7").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E7").Value = 'SOTENÄS'
$ws.Range("G7").Value = 0.5
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0

# Row 8: A 43053-2021
$ws.Range("A8").Value = 'A 43053-2021'
$ws.Range("B8").Value = 44431
$ws.Range("C8").Value = 46074
$ws.Range("D8").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E8").Value = 'SOTENÄS'
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0

# Row 9: A 43056-2021
$ws.Range("A9").Value = 'A 43056-2021'
$ws.Range("B9").Value = 44431
$ws.Range("C9").Value = 46074
$ws.Range("D9").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E9").Value = 'SOTENÄS'
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0

# Row 10: A 21487-2021
$ws.Range("A10").Value = 'A 21487-2021'
$ws.Range("B10").Value = 44320
$ws.Range("C10").Value = 46074
$ws.Range("D10").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E10").Value = 'SOTENÄS'
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0

# Row 11: A 73020-2021
$ws.Range("A11").Value = 'A 73020-2021'
$ws.Range("B11").Value = 44550
$ws.Range("C11").Value = 46074
$ws.Range("D11").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E11").Value = 'SOTENÄS'
$ws.Range("G11").Value = 2.7
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0

# Row 12: A 59231-2024
$ws.Range("A12").Value = 'A 59231-2024'
$ws.Range("B12").Value = 45637.58472222222
$ws.Range("C12").Value = 46074
$ws.Range("D12").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E12").Value = 'SOTENÄS'
$ws.Range("G12").Value = 1.3
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0

# Row 13: A 270-2025
$ws.Range("A13").Value = 'A 270-2025'
$ws.Range("B13").Value = 45660.48087962963
$ws.Range("C13").Value = 46074
$ws.Range("D13").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E13").Value = 'SOTENÄS'
$ws.Range("G13").Value = 8.9
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0

# Row 14: A 49633-2024
$ws.Range("A14").Value = 'A 49633-2024'
$ws.Range("B14").Value = 45596.59559027778
$ws.Range("C14").Value = 46074
$ws.Range("D14").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E14").Value = 'SOTENÄS'
$ws.Range("G14").Value = 0.8
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0

# Row 15: A 28416-2024
$ws.Range("A15").Value = 'A 28416-2024'
$ws.Range("B15").Value = 45477.62978009259
$ws.Range("C15").Value = 46074
$ws.Range("D15").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E15").Value = 'SOTENÄS'
$ws.Range("G15").Value = 0.2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0

# Row 16: A 4422-2024
$ws.Range("A16").Value = 'A 4422-2024'
$ws.Range("B16").Value = 45327.45375
$ws.Range("C16").Value = 46074
$ws.Range("D16").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E16").Value = 'SOTENÄS'
$ws.Range("G16").Value = 4.5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0

# Row 17: A 24212-2023
$ws.Range("A17").Value = 'A 24212-2023'
$ws.Range("B17").Value = 45076
$ws.Range("C17").Value = 46074
$ws.Range("D17").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E17").Value = 'SOTENÄS'
$ws.Range("G17").Value = 5.8
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0

# Row 18: A 50239-2022
$ws.Range("A18").Value = 'A 50239-2022'
$ws.Range("B18").Value = 44865
$ws.Range("C18").Value = 46074
$ws.Range("D18").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E18").Value = 'SOTENÄS'
$ws.Range("G18").Value = 13.2
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0

# Row 19: A 46579-2024
$ws.Range("A19").Value = 'A 46579-2024'
$ws.Range("B19").Value = 45582.75018518518
$ws.Range("C19").Value = 46074
$ws.Range("D19").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E19").Value = 'SOTENÄS'
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0

# Row 20: A 28409-2024
$ws.Range("A20").Value = 'A 28409-2024'
$ws.Range("B20").Value = 45477.62280092593
$ws.Range("C20").Value = 46074
$ws.Range("D20").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E20").Value = 'SOTENÄS'
$ws.Range("G20").Value = 0.4
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0

# Row 21: A 51434-2025
$ws.Range("A21").Value = 'A 51434-2025'
$ws.Range("B21").Value = 45949
$ws.Range("C21").Value = 46074
$ws.Range("D21").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E21").Value = 'SOTENÄS'
$ws.Range("G21").Value = 2.8
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0

# Row 22: A 52965-2025
$ws.Range("A22").Value = 'A 52965-2025'
$ws.Range("B22").Value = 45956
$ws.Range("C22").Value = 46074
$ws.Range("D22").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E22").Value = 'SOTENÄS'
$ws.Range("G22").Value = 0.6
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0

# Row 23: A 52960-2025
$ws.Range("A23").Value = 'A 52960-2025'
$ws.Range("B23").Value = 45956
$ws.Range("C23").Value = 46074
$ws.Range("D23").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E23").Value = 'SOTENÄS'
$ws.Range("G23").Value = 1.7
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0

# Row 24: A 52888-2025
$ws.Range("A24").Value = 'A 52888-2025'
$ws.Range("B24").Value = 45957.56943287037
$ws.Range("C24").Value = 46074
$ws.Range("D24").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E24").Value = 'SOTENÄS'
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0

# Row 25: A 57001-2025
$ws.Range("A25").Value = 'A 57001-2025'
$ws.Range("B25").Value = 45977
$ws.Range("C25").Value = 46074
$ws.Range("D25").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E25").Value = 'SOTENÄS'
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0

# Row 26: A 4780-2022
$ws.Range("A26").Value = 'A 4780-2022'
$ws.Range("B26").Value = 44592.62657407407
$ws.Range("C26").Value = 46074
$ws.Range("D26").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E26").Value = 'SOTENÄS'
$ws.Range("G26").Value = 0.9
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0

# Row 27: A 17492-2024
$ws.Range("A27").Value = 'A 17492-2024'
$ws.Range("B27").Value = 45415.50709490741
$ws.Range("C27").Value = 46074
$ws.Range("D27").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E27").Value = 'SOTENÄS'
$ws.Range("G27").Value = 7.7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0

# Row 28: A 46587-2024
$ws.Range("A28").Value = 'A 46587-2024'
$ws.Range("B28").Value = 45582.76763888889
$ws.Range("C28").Value = 46074
$ws.Range("D28").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E28").Value = 'SOTENÄS'
$ws.Range("G28").Value = 2.8
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0

# Row 29: A 46588-2024
$ws.Range("A29").Value = 'A 46588-2024'
$ws.Range("B29").Value = 45582.77137731481
$ws.Range("C29").Value = 46074
$ws.Range("D29").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E29").Value = 'SOTENÄS'
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0

# Row 30: A 28418-2024
$ws.Range("A30").Value = 'A 28418-2024'
$ws.Range("B30").Value = 45477.62978009259
$ws.Range("C30").Value = 46074
$ws.Range("D30").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E30").Value = 'SOTENÄS'
$ws.Range("G30").Value = 0.2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0

# Row 31: A 49634-2024
$ws.Range("A31").Value = 'A 49634-2024'
$ws.Range("B31").Value = 45596.59591435185
$ws.Range("C31").Value = 46074
$ws.Range("D31").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E31").Value = 'SOTENÄS'
$ws.Range("G31").Value = 0.7
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0

# Row 32: A 21421-2021
$ws.Range("A32").Value = 'A 21421-2021'
$ws.Range("B32").Value = 44316
$ws.Range("C32").Value = 46074
$ws.Range("D32").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E32").Value = 'SOTENÄS'
$ws.Range("G32").Value = 0.6
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0

# Row 33: A 17491-2024
$ws.Range("A33").Value = 'A 17491-2024'
$ws.Range("B33").Value = 45415.50266203703
$ws.Range("C33").Value = 46074
$ws.Range("D33").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E33").Value = 'SOTENÄS'
$ws.Range("G33").Value = 6.2
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0

# Row 34: A 50864-2022
$ws.Range("A34").Value = 'A 50864-2022'
$ws.Range("B34").Value = 44867.56143518518
$ws.Range("C34").Value = 46074
$ws.Range("D34").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E34").Value = 'SOTENÄS'
$ws.Range("G34").Value = 3.3
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0

# Row 35: A 23503-2025
$ws.Range("A35").Value = 'A 23503-2025'
$ws.Range("B35").Value = 45795
$ws.Range("C35").Value = 46074
$ws.Range("D35").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E35").Value = 'SOTENÄS'
$ws.Range("G35").Value = 14.1
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0

# Row 36: A 24086-2025
$ws.Range("A36").Value = 'A 24086-2025'
$ws.Range("B36").Value = 45795
$ws.Range("C36").Value = 46074
$ws.Range("D36").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E36").Value = 'SOTENÄS'
$ws.Range("G36").Value = 0.7
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0

